# Addressbook.xlsx — add a new "Employee Email" field/column between
# "PhoneNumber" and "ID" (mirrors the commit: "Add new identical field
# for Addressbook and english form: Employee Email").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting the old ID/Status columns to D/E.
$ws.Columns("C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "Employee Email"

# Restore/update the worksheet selection the author ended up with.
[void]$ws.Range("A2:B3").Select()
